$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -81.6599
$ws.Range("B2").Value = -81.5538

$ws.Range("A3").Value = 29.5627
$ws.Range("B3").Value = 29.6557

$ws.Range("A4").Value = -80.919
$ws.Range("B4").Value = -81.0261

$ws.Range("A5").Value = 30.2073
$ws.Range("B5").Value = 30.1148
